$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (289:290), pushing the
# existing rows 289-298 down to 291-300.
$ws.Rows("289:290").Insert()

# New row 289: week of 2021-11-09, "1a amarillo"
$ws.Range("A289").Value = 11
$ws.Range("B289").Value = "Vega Monumental Concepción"
$ws.Range("C289").Value = "Bíobío"
$ws.Range("D289").Value = 44509
$ws.Range("E289").Value = 8
$ws.Range("F289").Value = "Fruta"
$ws.Range("G289").Value = 100102
$ws.Range("H289").Value = "Cítricos"
$ws.Range("I289").Value = 100102003
$ws.Range("J289").Value = "Limón"
$ws.Range("K289").Value = "Sin especificar"
$ws.Range("L289").Value = "1a amarillo"
$ws.Range("M289").Value = 510
$ws.Range("N289").Value = 8000
$ws.Range("O289").Value = 8500
$ws.Range("P289").Value = 8255
$ws.Range("Q289").Value = "`$/malla 16 kilos"
$ws.Range("R289").Value = "Región de O'Higgins"
$ws.Range("S289").Value = 516
$ws.Range("T289").Value = 16

# New row 290: week of 2021-11-09, "2a amarillo"
$ws.Range("A290").Value = 11
$ws.Range("B290").Value = "Vega Monumental Concepción"
$ws.Range("C290").Value = "Bíobío"
$ws.Range("D290").Value = 44509
$ws.Range("E290").Value = 8
$ws.Range("F290").Value = "Fruta"
$ws.Range("G290").Value = 100102
$ws.Range("H290").Value = "Cítricos"
$ws.Range("I290").Value = 100102003
$ws.Range("J290").Value = "Limón"
$ws.Range("K290").Value = "Sin especificar"
$ws.Range("L290").Value = "2a amarillo"
$ws.Range("M290").Value = 200
$ws.Range("N290").Value = 6500
$ws.Range("O290").Value = 6500
$ws.Range("P290").Value = 6500
$ws.Range("Q290").Value = "`$/malla 16 kilos"
$ws.Range("R290").Value = "Región de O'Higgins"
$ws.Range("S290").Value = 406
$ws.Range("T290").Value = 16
